$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Update the cached "datetimeFigureOut" date field text (1/22/2020 ->
#    1/28/2020) on the slide master and every slide layout's Date
#    Placeholder shape.
# ---------------------------------------------------------------------------

function Update-DateField($shapes, $newText) {
    for ($j = 1; $j -le $shapes.Count; $j++) {
        $shp = $shapes.Item($j)
        if ($shp.Name -like "Date Placeholder*") {
            $tr = $shp.TextFrame.TextRange
            $para = $tr.Paragraphs(1, 1)
            $run = $para.Runs(1, 1)
            $run.Text = "1/28/2020"
        }
    }
}

$master = $p.SlideMaster
Update-DateField $master.Shapes "1/28/2020"

$layouts = $master.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    $layout = $layouts.Item($i)
    Update-DateField $layout.Shapes "1/28/2020"
}

# ---------------------------------------------------------------------------
# 2) Slide 9 ("Platform"), Content Placeholder 2, bullet "Chance to hire
#    contractor programmer" -> "... with lower cost and more support"
# ---------------------------------------------------------------------------

$s = $p.Slides.Item(9)
$shp = $s.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange
$para = $tr.Paragraphs(7, 1)
$run = $para.Runs(1, 1)
$run.Text = "Chance to hire contractor programmer with lower cost and more support"
